# Apply "carl's newest scraper" layout change:
#  - Drop the extra/legacy columns that the old scraper emitted
#    (Norm-Typ, Ritningsnummer, Position, Beteckning, Kompletterande
#    Information ovrigt, Ref annan, Historiskt Varumarke, Historiskt
#    inkopsreferens, Forpackning and the empty spacer column) so that
#    only Varumarke..RSK-nummer remain, packed into columns A:I.
#  - Give the new "SSG-notering" column (now G) a wider, explicit width.
#  - Turn off the AutoFilter / sort state that referenced the old A:S range.
#  - Re-point the hidden _FilterDatabase defined name at the new A:I range.
#  - Reset the view so it is scrolled back to A1 and the active selection
#    follows the (now) E-nummer column H instead of the old column Q.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the obsolete columns. Go from right to left so earlier deletions
# don't shift the letters of columns we still need to delete.
$obsoleteColumns = @("Q", "O", "M", "L", "K", "J", "I", "H", "G", "A")
foreach ($col in $obsoleteColumns) {
    $ws.Columns($col).Delete()
}

# The surviving "SSG-notering" column (originally P, now G) gets an
# explicit custom width.
$ws.Columns(7).ColumnWidth = 28.65

# Drop the old AutoFilter / sortState (it referenced the stale A1:S21217
# range and is no longer wanted on the trimmed table).
$ws.AutoFilterMode = $false

# Re-point the hidden filter-database defined name at the new, narrower
# range now that columns have been removed.
foreach ($n in $wb.Names) {
    if ($n.Name -match "_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$I`$21217"
    }
}

# Reset the window scroll position back to the top-left (A1) and move the
# active column selection from the old column Q to the new column H.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
[void]$ws.Range("H1:H1048576").Select()
